$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Instructions")
$ws2.Range("A6").Value = '4. Save as the excel file (only the "Fill out this form" sheet)  as "Text (tab-delimited) (*.txt)" or "CSV (comma-delimited) (*.csv)"'
$ws2.Activate()
$ws2.Range("A6").Select()
